$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''329.10'
$ws.Range("E2").Value = '''0.34%'
$ws.Range("G2").Value = '''15'
$ws.Range("D3").Value = '''44.33'
$ws.Range("E3").Value = '''0.82%'
$ws.Range("G3").Value = '''15'
$ws.Range("D4").Value = '''5.483'
$ws.Range("E4").Value = '''-1.43%'
$ws.Range("G4").Value = '''15'
$ws.Range("D5").Value = '''0.08075'
$ws.Range("E5").Value = '''0.35%'
$ws.Range("G5").Value = '''15'
$ws.Range("D6").Value = '''2.047'
$ws.Range("E6").Value = '''7.93%'
$ws.Range("G6").Value = '''15'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = '''0.9533'
$ws.Range("E7").Value = '''0.84%'
$ws.Range("G7").Value = '''15'
$ws.Range("B8").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C8").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D8").Value = '''0.1132'
$ws.Range("E8").Value = '''-3.65%'
$ws.Range("G8").Value = '''15'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '''0.1890'
$ws.Range("E9").Value = '''2.73%'
$ws.Range("G9").Value = '''15'
$ws.Range("B10").Value = 'MCDex'
$ws.Range("C10").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D10").Value = '''10.20'
$ws.Range("E10").Value = '''-6.26%'
$ws.Range("G10").Value = '''15'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.09969'
$ws.Range("E11").Value = '''2.99%'
$ws.Range("G11").Value = '''15'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.04798'
$ws.Range("E12").Value = '''10.00%'
$ws.Range("G12").Value = '''15'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.1065'
$ws.Range("E13").Value = '''-0.24%'
$ws.Range("G13").Value = '''15'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001266'
$ws.Range("E14").Value = '''-1.50%'
$ws.Range("G14").Value = '''15'
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").Value = '''0.04089'
$ws.Range("E15").Value = '''-2.13%'
$ws.Range("G15").Value = '''15'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005904'
$ws.Range("E16").Value = '''-1.46%'
$ws.Range("G16").Value = '''15'
$ws.Range("B17").Value = 'HotbitToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D17").Value = '''0.004359'
$ws.Range("E17").Value = '''1.75%'
$ws.Range("G17").Value = '''15'
$ws.Range("D18").Value = '''3.377'
$ws.Range("E18").Value = '''-0.90%'
$ws.Range("G18").Value = '''15'
$ws.Range("D19").Value = '''4.420'
$ws.Range("E19").Value = '''3.51%'
$ws.Range("G19").Value = '''15'
$ws.Range("B20").Value = 'BTSEToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D20").Value = '''2.580'
$ws.Range("E20").Value = '''1.62%'
$ws.Range("G20").Value = '''15'
$ws.Range("B21").Value = 'BitpandaEcosystemToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D21").Value = '''0.3379'
$ws.Range("E21").Value = '''-1.99%'
$ws.Range("G21").Value = '''15'
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").Value = '''0.1401'
$ws.Range("E22").Value = '''0.92%'
$ws.Range("G22").Value = '''15'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '''0.2579'
$ws.Range("E23").Value = '''2.96%'
$ws.Range("G23").Value = '''15'
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").Value = '''0.001307'
$ws.Range("E24").Value = '''4.83%'
$ws.Range("G24").Value = '''15'
$ws.Range("D25").Value = '''0.0001251'
$ws.Range("E25").Value = '''-0.75%'
$ws.Range("G25").Value = '''15'
$ws.Range("E26").Value = '''-6.20%'
$ws.Range("G26").Value = '''15'
$ws.Range("G27").Value = '''15'
$ws.Range("G28").Value = '''15'
$ws.Range("G29").Value = '''15'
$ws.Range("G30").Value = '''15'
$ws.Range("G31").Value = '''15'
$ws.Range("G32").Value = '''15'
$ws.Range("G33").Value = '''15'
$ws.Range("G34").Value = '''15'
$ws.Range("G35").Value = '''15'
$ws.Range("G36").Value = '''15'
$ws.Range("G37").Value = '''15'
$ws.Range("D38").Value = '''0.02591'
$ws.Range("E38").Value = '''-2.08%'
$ws.Range("G38").Value = '''15'
$ws.Range("D39").Value = '''0.05684'
$ws.Range("E39").Value = '''3.22%'
$ws.Range("G39").Value = '''15'
$ws.Range("D40").Value = '''0.007573'
$ws.Range("E40").Value = '''-0.25%'
$ws.Range("G40").Value = '''15'
$ws.Range("D41").Value = '''0.1400'
$ws.Range("E41").Value = '''0.40%'
$ws.Range("G41").Value = '''15'
$ws.Range("D42").Value = '''0.007344'
$ws.Range("E42").Value = '''-7.36%'
$ws.Range("G42").Value = '''15'
$ws.Range("D43").Value = '''0.002009'
$ws.Range("E43").Value = '''0.42%'
$ws.Range("G43").Value = '''15'
$ws.Range("D44").Value = '''0.009053'
$ws.Range("E44").Value = '''2.47%'
$ws.Range("G44").Value = '''15'
$ws.Range("D45").Value = '''0.00007016'
$ws.Range("E45").Value = '''1.84%'
$ws.Range("G45").Value = '''15'
$ws.Range("E46").Value = '''0.06%'
$ws.Range("G46").Value = '''15'
$ws.Range("E47").Value = '''-0.09%'
$ws.Range("G47").Value = '''15'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").Value = '''0.003509'
$ws.Range("E48").Value = '''0.24%'
$ws.Range("G48").Value = '''15'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").Value = '''0.003502'
$ws.Range("E49").Value = '''54.18%'
$ws.Range("G49").Value = '''15'
$ws.Range("E50").Value = '''0.06%'
$ws.Range("G50").Value = '''15'
$ws.Range("E51").Value = '''0.06%'
$ws.Range("G51").Value = '''15'
